# Renumber the "Decision Tree" chapter from 12 -> 13 throughout the deck.
# Every slide from 1 to 22 has its title placeholder (the first shape on
# the slide) starting with "12" (e.g. "12 Introduction to Decision Tree",
# "12.1 Decision Tree Algorithm", "12.2 Build Tree", ...). Bump the leading
# chapter number to 13 while leaving the rest of the title untouched.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $titleShape = $slide.Shapes.Item(1)

    if ($titleShape.HasTextFrame) {
        $tr = $titleShape.TextFrame.TextRange
        $text = $tr.Text

        if ($text -match '^12(\D|$)') {
            $newText = '13' + $text.Substring(2)
            $tr.Text = $newText
        }
    }
}
